$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("Aktif" shifts from F to G),
# making room for the new "Tambahan Solar" (solar cost) column.
$ws.Columns("F:F").Insert()

# Populate the header of the newly inserted column.
$ws.Range("F3").Value = "Tambahan Solar"

# Match the author's final selection/cursor position.
$ws.Range("F3").Select()
